# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 359
$firstRow = 2

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45188
